# Applies the cryptos-list price/volume refresh described by the commit diff.
# Numeric-looking strings (e.g. "2.60", "415.47") must stay TEXT cells (as in
# the source data), so we write them with a leading apostrophe (forces Excel
# to treat the input as text instead of auto-converting to a number) and then
# reset the cell style back to "Normal" so no stray number-format/quote-prefix
# styling is left behind on the cell.
function Set-TextValue {
    param($ws, $addr, $val)
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws "D2" "67.280.84"
Set-TextValue $ws "E2" "  +7.13%  "

Set-TextValue $ws "D3" "3.585.45"
Set-TextValue $ws "E3" "  +3.40%  "

Set-TextValue $ws "E4" "  +0.14%  "

Set-TextValue $ws "D5" "415.47"
Set-TextValue $ws "E5" "  +0.22%  "

Set-TextValue $ws "D6" "129.55"
Set-TextValue $ws "E6" "  -0.24%  "

Set-TextValue $ws "D7" "0.650"
Set-TextValue $ws "E7" "  +3.67%  "

Set-TextValue $ws "D8" "3.579.56"
Set-TextValue $ws "E8" "  +3.41%  "

Set-TextValue $ws "D9" "0.999"
Set-TextValue $ws "E9" "  -0.07%  "

Set-TextValue $ws "D10" "0.770"
Set-TextValue $ws "E10" "  +6.26%  "

Set-TextValue $ws "D11" "0.177"
Set-TextValue $ws "E11" "  +15.87%  "

Set-TextValue $ws "D12" "0.0000331"
Set-TextValue $ws "E12" "  +44.24%  "

Set-TextValue $ws "D13" "42.22"
Set-TextValue $ws "E13" "  -0.59%  "

Set-TextValue $ws "D14" "9.89"
Set-TextValue $ws "E14" "  +2.09%  "

Set-TextValue $ws "D15" "4.157.84"
Set-TextValue $ws "E15" "  +3.41%  "

Set-TextValue $ws "E16" "  -0.19%  "

Set-TextValue $ws "D17" "20.36"
Set-TextValue $ws "E17" "  -0.60%  "

Set-TextValue $ws "D18" "3.569.92"
Set-TextValue $ws "E18" "  +3.12%  "

Set-TextValue $ws "D19" "1.13"

Set-TextValue $ws "D20" "67.340.72"
Set-TextValue $ws "E20" "  +7.33%  "

Set-TextValue $ws "E21" "  -2.72%  "

Set-TextValue $ws "D22" "451.74"
Set-TextValue $ws "E22" "  -1.71%  "

Set-TextValue $ws "D23" "88.69"
Set-TextValue $ws "E23" "  -1.79%  "

Set-TextValue $ws "E24" "  -4.13%  "

Set-TextValue $ws "D25" "13.12"
Set-TextValue $ws "E25" "  -0.67%  "

Set-TextValue $ws "E26" "  +1.17%  "

Set-TextValue $ws "E27" "  -6.64%  "

Set-TextValue $ws "D28" "34.85"
Set-TextValue $ws "E28" "  +4.67%  "

Set-TextValue $ws "D29" "4.88"
Set-TextValue $ws "E29" "  +1.88%  "

Set-TextValue $ws "E30" "  +4.16%  "

Set-TextValue $ws "D31" "12.32"
Set-TextValue $ws "E31" "  +2.22%  "

Set-TextValue $ws "E32" "  +4.66%  "

Set-TextValue $ws "D33" "7.39"
Set-TextValue $ws "E33" "  -2.13%  "

Set-TextValue $ws "D34" "0.161"
Set-TextValue $ws "E34" "  -4.03%  "

Set-TextValue $ws "D35" "40.73"
Set-TextValue $ws "E35" "  -0.12%  "

Set-TextValue $ws "E36" "  -0.04%  "

Set-TextValue $ws "D37" "56.70"
Set-TextValue $ws "E37" "  -2.76%  "

Set-TextValue $ws "D38" "0.0492"
Set-TextValue $ws "E38" "  +0.67%  "

Set-TextValue $ws "D39" "0.0₃0748"
Set-TextValue $ws "E39" "  +33.23%  "

Set-TextValue $ws "D41" "0.999"
Set-TextValue $ws "E41" "  -0.07%  "

Set-TextValue $ws "D42" "3.04"
Set-TextValue $ws "E42" "  -0.95%  "

Set-TextValue $ws "E43" "  +1.44%  "

Set-TextValue $ws "D44" "148.99"
Set-TextValue $ws "E44" "  -0.57%  "

Set-TextValue $ws "E45" "  -2.30%  "

Set-TextValue $ws "D46" "3.25"
Set-TextValue $ws "E46" "  -1.90%  "

Set-TextValue $ws "D47" "4.30"
Set-TextValue $ws "E47" "  -2.49%  "

Set-TextValue $ws "D48" "1.97"
Set-TextValue $ws "E48" "  -4.01%  "

Set-TextValue $ws "D49" "2.34"
Set-TextValue $ws "E49" "  -1.35%  "

$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws "D50" "2.60"
Set-TextValue $ws "E50" "  +11.33%  "

$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws "D51" "114.72"
Set-TextValue $ws "E51" "  +5.46%  "
